# Horarios actualizados Linea 141 - 767
# Updates the scraped-schedule workbook: refreshes the "last updated" timestamp
# (03:57:17 -> 04:24:09), updates row counts / minute-remaining values for rows
# whose Hora_Scrap was still the old timestamp, and appends newly scraped rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value2 = "Última actualización: 04:24:09"
$ws1.Cells.Item(3, 1).Value2 = "Total filas: 26"

# Existing rows whose Hora_Scrap refreshed to the new scrape time, with
# updated Minutos countdown.
$ws1.Cells.Item(19, 1).Value2 = "04:24:09"
$ws1.Cells.Item(19, 4).Value2 = 29

$ws1.Cells.Item(21, 1).Value2 = "04:24:09"
$ws1.Cells.Item(21, 4).Value2 = 52

$ws1.Cells.Item(22, 1).Value2 = "04:24:09"
$ws1.Cells.Item(22, 4).Value2 = 58

$ws1.Cells.Item(23, 1).Value2 = "04:24:09"
$ws1.Cells.Item(23, 4).Value2 = 71

$ws1.Cells.Item(25, 1).Value2 = "04:24:09"
$ws1.Cells.Item(25, 4).Value2 = 82

# Newly scraped rows appended at the bottom.
$ws1.Cells.Item(26, 1).Value2 = "04:24:09"
$ws1.Cells.Item(26, 2).Value2 = "05:48"
$ws1.Cells.Item(26, 3).Value2 = "215A_EL PATO"
$ws1.Cells.Item(26, 4).Value2 = 84
$ws1.Cells.Item(26, 5).Value2 = "LP1912"

$ws1.Cells.Item(27, 1).Value2 = "04:24:09"
$ws1.Cells.Item(27, 2).Value2 = "05:54"
$ws1.Cells.Item(27, 3).Value2 = "10_OLMOS"
$ws1.Cells.Item(27, 4).Value2 = 90
$ws1.Cells.Item(27, 5).Value2 = "LP1912"

$ws1.Cells.Item(28, 1).Value2 = "04:24:09"
$ws1.Cells.Item(28, 2).Value2 = "06:09"
$ws1.Cells.Item(28, 3).Value2 = "16_SANTA ANA"
$ws1.Cells.Item(28, 4).Value2 = 105
$ws1.Cells.Item(28, 5).Value2 = "LP1912"

$ws1.Cells.Item(29, 1).Value2 = "04:24:09"
$ws1.Cells.Item(29, 2).Value2 = "06:11"
$ws1.Cells.Item(29, 3).Value2 = "215A_EL PATO"
$ws1.Cells.Item(29, 4).Value2 = 107
$ws1.Cells.Item(29, 5).Value2 = "LP1912"

$ws1.Cells.Item(30, 1).Value2 = "04:24:09"
$ws1.Cells.Item(30, 2).Value2 = "06:14"
$ws1.Cells.Item(30, 3).Value2 = "225_HARAS DEL SUR"
$ws1.Cells.Item(30, 4).Value2 = 110
$ws1.Cells.Item(30, 5).Value2 = "LP1912"

$ws1.Cells.Item(31, 1).Value2 = "04:24:09"
$ws1.Cells.Item(31, 2).Value2 = "06:21"
$ws1.Cells.Item(31, 3).Value2 = "26_HERNANDEZ"
$ws1.Cells.Item(31, 4).Value2 = 117
$ws1.Cells.Item(31, 5).Value2 = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value2 = "Última actualización: 04:24:09"
$ws2.Cells.Item(3, 1).Value2 = "Total filas: 9"

$ws2.Cells.Item(12, 1).Value2 = "04:24:09"
$ws2.Cells.Item(12, 4).Value2 = 71

$ws2.Cells.Item(13, 1).Value2 = "04:24:09"
$ws2.Cells.Item(13, 2).Value2 = "05:48"
$ws2.Cells.Item(13, 3).Value2 = "215A_EL PATO"
$ws2.Cells.Item(13, 4).Value2 = 84
$ws2.Cells.Item(13, 5).Value2 = "LP1912"

$ws2.Cells.Item(14, 1).Value2 = "04:24:09"
$ws2.Cells.Item(14, 2).Value2 = "06:11"
$ws2.Cells.Item(14, 3).Value2 = "215A_EL PATO"
$ws2.Cells.Item(14, 4).Value2 = 107
$ws2.Cells.Item(14, 5).Value2 = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value2 = "Última actualización: 04:24:09"

$ws3.Cells.Item(7, 1).Value2 = "04:24:09"
$ws3.Cells.Item(7, 4).Value2 = 80
